$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 45815.39107167824

$ws.Range("A16").Value = 45816.39137637648
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").Value = "EVOWHEY PROTEIN"
$ws.Range("C16").Value = "2Kg"
$ws.Range("D16").Value = "34,90€"
